# Scale the "value" column (D) by 10000 for all data rows (2-33).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 33; $row++) {
    $cell = $ws.Cells.Item($row, 4)
    $cell.Value = $cell.Value2 * 10000
}
